$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set sum values for B column (previously empty numeric cells)
$ws.Range("B3").Value = 65
$ws.Range("B6").Value = 44
$ws.Range("B8").Value = 55
$ws.Range("B10").Value = 15

# Update model names in column A
$ws.Range("A8").Value = "Statyw drewniany"
$ws.Range("A10").Value = "Statyw metalowy"
